# Applies the "Refined metadata to be additional tab" edit:
#   1. Refreshes the panel-query timestamps recorded in the "data" sheet
#      (F2:F35) to the time of the later re-run.
#   2. Adds a new "metadata" worksheet (placed after "data") describing the
#      PanelApp query that produced the data, mirroring the header/row
#      formatting used on the "data" sheet.

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("data")

# --- 1. Refresh panel_query_time / time_taken column on "data" ------------
$dataSheet.Range("F2").Value = "2021-10-05 14:19:43.951075"
$dataSheet.Range("F3").Value = "2021-10-05 14:19:43.951083"
$dataSheet.Range("F4").Value = "2021-10-05 14:19:43.951086"
$dataSheet.Range("F5").Value = "2021-10-05 14:19:43.951089"
$dataSheet.Range("F6").Value = "2021-10-05 14:19:43.951092"
$dataSheet.Range("F7").Value = "2021-10-05 14:19:43.951094"
$dataSheet.Range("F8").Value = "2021-10-05 14:19:43.951097"
$dataSheet.Range("F9").Value = "2021-10-05 14:19:43.951099"
$dataSheet.Range("F10").Value = "2021-10-05 14:19:43.951102"
$dataSheet.Range("F11").Value = "2021-10-05 14:19:43.951105"
$dataSheet.Range("F12").Value = "2021-10-05 14:19:43.951107"
$dataSheet.Range("F13").Value = "2021-10-05 14:19:43.951110"
$dataSheet.Range("F14").Value = "2021-10-05 14:19:43.951112"
$dataSheet.Range("F15").Value = "2021-10-05 14:19:43.951114"
$dataSheet.Range("F16").Value = "2021-10-05 14:19:43.951117"
$dataSheet.Range("F17").Value = "2021-10-05 14:19:43.951119"
$dataSheet.Range("F18").Value = "2021-10-05 14:19:43.951122"
$dataSheet.Range("F19").Value = "2021-10-05 14:19:43.951125"
$dataSheet.Range("F20").Value = "2021-10-05 14:19:43.951127"
$dataSheet.Range("F21").Value = "2021-10-05 14:19:43.951130"
$dataSheet.Range("F22").Value = "2021-10-05 14:19:43.951132"
$dataSheet.Range("F23").Value = "2021-10-05 14:19:43.951135"
$dataSheet.Range("F24").Value = "2021-10-05 14:19:43.951137"
$dataSheet.Range("F25").Value = "2021-10-05 14:19:43.951140"
$dataSheet.Range("F26").Value = "2021-10-05 14:19:43.951142"
$dataSheet.Range("F27").Value = "2021-10-05 14:19:43.951145"
$dataSheet.Range("F28").Value = "2021-10-05 14:19:43.951148"
$dataSheet.Range("F29").Value = "2021-10-05 14:19:43.951150"
$dataSheet.Range("F30").Value = "2021-10-05 14:19:43.951153"
$dataSheet.Range("F31").Value = "2021-10-05 14:19:43.951155"
$dataSheet.Range("F32").Value = "2021-10-05 14:19:43.951158"
$dataSheet.Range("F33").Value = "2021-10-05 14:19:43.951161"
$dataSheet.Range("F34").Value = "2021-10-05 14:19:43.951164"
$dataSheet.Range("F35").Value = "2021-10-05 14:19:43.951166"

# --- 2. Add the "metadata" worksheet, right after "data" ------------------
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$newSheet.Name = "metadata"

$newSheet.Range("B1").Value = "data_name"
$newSheet.Range("C1").Value = "data_id"
$newSheet.Range("D1").Value = "data_version"
$newSheet.Range("E1").Value = "data_version_created"
$newSheet.Range("F1").Value = "panel_query_time"
$newSheet.Range("G1").Value = "panel_get_request"

$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").Value = "Congenital hypothyroidism"
$newSheet.Range("C2").Value = 31
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "2.5"
$newSheet.Range("D2").Style = "Normal"
$newSheet.Range("E2").Value = "2021-08-24T12:30:05.851921Z"
$newSheet.Range("F2").Value = "2021-10-05 14:19:43.947296"
$newSheet.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/31/?format=json"

# Match the "data" sheet's header styling (bold, centred, thin border)
$dataSheet.Range("B1").Copy()
$newSheet.Range("B1:G1").PasteSpecial(-4122)
$dataSheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)
